$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: month number (C9) drives the A3 "Wykaz Czynnosci" title formula ---
$ws.Range("C9").Value = 5

# --- Week-number column (A) for each week block ---
$ws.Range("A11").Value = 18
$ws.Range("A18").Value = 19
$ws.Range("A25").Value = 20
$ws.Range("A32").Value = 21
$ws.Range("A39").Value = 22

# --- Anchor date shifts forward 13 weeks (91 days); dependent B-column formulas cascade ---
$ws.Range("B11").Value = 42856

# --- Row 46 total ---
$ws.Range("C46").Value = 21

# --- Unhide weekly detail rows that become visible again (11, 12, 41) ---
$ws.Rows(11).Hidden = $false
$ws.Rows(12).Hidden = $false
$ws.Rows(41).Hidden = $false

# --- Per-day Miejsce pracy (D=obecnosc, F=opis dzialalnosci, E=miejsce pracy) ---
# Row 12
$ws.Range("D12").Value = 1
$ws.Range("F12").Value = 'Bonikowo. Weryfikacja wniosku o zatwierdzenie podwykonawcy - Endcor'
$ws.Range("E12").Value = 'Bonikowo'
# Row 13
$ws.Range("D13").Value = ""
$ws.Range("F13").Value = ""
# Row 14
$ws.Range("F14").Value = 'Bonikowo. Raport BHP za miesiąc kwiecień'
$ws.Range("E14").Value = 'Bonikowo'
# Row 15
$ws.Range("F15").Value = 'Bonikowo. Wizyta na budowie w czasie odhumusowania, protokół z kontroli inspekcyjnej bhp z dnia 05.05.2017'
$ws.Range("E15").Value = 'Bonikowo'
# Row 18
$ws.Range("F18").Value = 'Bonikowo. Weryfikacja wniosków o zatwierdzenie podwykonawcy (Ekoinvest, Wprinż, Nodic)'
$ws.Range("E18").Value = 'Bonikowo'
# Row 19
$ws.Range("F19").Value = 'Bonikowo. Weryfikacja wniosków o zatwierdzenie podwykonawcy (Metkor, Ornia, Rafaco, Saniment)'
$ws.Range("E19").Value = 'Bonikowo'
# Row 20
$ws.Range("F20").Value = 'Bonikowo. Narada koordynacyjna. Wizyta na budowie, protokół z kontroli bhp z dnia 10.05.2017'
$ws.Range("E20").Value = 'Bonikowo'
# Row 21
$ws.Range("F21").Value = 'Bonikowo. Udział w tygodniu bezpieczeństaw - pokaz POPŻ. Weryfikacja IT 5.4.23'
$ws.Range("E21").Value = 'Bonikowo'
# Row 22
$ws.Range("F22").Value = 'Bonikowo. Weryfikacja protokołu dopuszczenia Budinż.'
$ws.Range("E22").Value = 'Bonikowo'
# Row 25
$ws.Range("F25").Value = 'Bonikowo. Wizyta na budowie, protokół z kontroli bhp z dnia 15.05.17'
$ws.Range("E25").Value = 'Bonikowo'
# Row 26
$ws.Range("F26").Value = 'Bonikowo.Weryfikacja wniosków o zatwierdzenie podwykonawcy (Wunderlich, Izostal)'
$ws.Range("E26").Value = 'Bonikowo'
# Row 27
$ws.Range("F27").Value = 'Bonikowo. Wizyta na budowie. Protokół z kontroli bhp z dnia 17.05.17'
$ws.Range("E27").Value = 'Bonikowo'
# Row 28
$ws.Range("F28").Value = 'Bonikowo. Wizyta na budowie - przewierty próbne'
$ws.Range("E28").Value = 'Bonikowo'
# Row 29
$ws.Range("F29").Value = 'Bonikowo. Weryfikacja wniosków o zatwierdzenie podwykonawcy (Sznajder, Dzwigmar, ATS)'
$ws.Range("E29").Value = 'Bonikowo'
# Row 32
$ws.Range("D32").Value = 1
$ws.Range("F32").Value = 'Bonikowo. Weryfikacja wniosku o zatwierdzenie podwykonawcy AHAK, RENOMA. Sprawdzenie dokumntów bhp pracowników fizycznych KVV. Wizyta na budowie'
$ws.Range("E32").Value = 'Bonikowo'
# Row 33
$ws.Range("D33").Value = 1
$ws.Range("F33").Value = 'Bonikowo. Weryfikacja wniosku o zatwierdzenie podwykonawcy - Metkor. Protokół dopuszczenia - Sznajder'
$ws.Range("E33").Value = 'Bonikowo'
# Row 34
$ws.Range("D34").Value = 1
$ws.Range("F34").Value = 'Bonikowo. Rada budowy. Wizyta na budowie.'
$ws.Range("E34").Value = 'Bonikowo'
# Row 35
$ws.Range("D35").Value = 1
$ws.Range("F35").Value = 'Bonikowo. Sprawdzenie uprawnień seposkich pracowników fizycznych KVV.'
$ws.Range("E35").Value = 'Bonikowo'
# Row 36
$ws.Range("D36").Value = 1
$ws.Range("F36").Value = 'Bonikowo. Sprawdzenie zgodności kadry nadzorującej zgodnie ze schematem organizacyjnym. Sprawdzenie zgodności umowy Orina.'
$ws.Range("E36").Value = 'Bonikowo'
# Row 39
$ws.Range("D39").Value = 1
$ws.Range("F39").Value = 'Bonikowo. Wizyta na budowie - protokół z kontroli bhp z dnia 29.05.17'
$ws.Range("E39").Value = 'Bonikowo'
# Row 40
$ws.Range("D40").Value = 1
$ws.Range("F40").Value = 'Bonikowo. Weryfikacja protokołu dopuszczenia Budinż (dodatkowe osoby + sprzęt)'
$ws.Range("E40").Value = 'Bonikowo'
# Row 41
$ws.Range("D41").Value = 1
$ws.Range("F41").Value = 'Bonikowo. Narada Koordynacyjna. Raport BHP nr 4 K121-ILFWs-RM-0007'
$ws.Range("E41").Value = 'Bonikowo'

# --- Restore natural row height/remove custom-height artifact for newly unhidden rows ---
$ws.Rows(11).AutoFit()
$ws.Rows(12).AutoFit()
$ws.Rows(41).AutoFit()
